$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1, J1), reusing H1's style (border/bold/centered)
# by copying H1 onto each target cell before overwriting its value.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-20.
$iValues = @(5, 1, 3, 3, 6, 10, 1, 1, 2, 1, 4, 1, 1, 1, 1, 1, 1, 9, 9)
$jValues = @(6, 2, 5, 6, 8, 10, 4, 4, 5, 3, 5, 4, 4, 6, 6, 4, 2, 9, 9)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
